$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.627.59"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.677.98"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "595.30"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "174.89"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").Value = "2.678.33"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "4.98"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "3.171.52"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "71.570.12"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "26.04"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "2.643.76"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "12.02"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("D20").Value = "7.97"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "368.03"
$ws.Range("E21").Value = "  -3.56%  "
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "71.61"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "4.28"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").Value = "9.81"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "2.817.74"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "502.75"
$ws.Range("E32").Value = "  -8.47%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "162.67"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "19.36"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").Value = "1.36"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "0.108"
$ws.Range("E40").Value = "  -5.95%  "
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "4.97"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "155.74"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "39.17"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").Value = "3.69"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").Value = "0.0762"
$ws.Range("E51").Value = "  +0.47%  "
